$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells below hold plain-text numeric/percentage strings (t="inlineStr" in the
# source sheet). Force Text number format first so Excel does not coerce the
# assigned strings into numeric/percentage values, which would lose the exact
# textual representation (e.g. trailing zeros, literal "%" sign).
$updates = @{
    "D2" = "307.23"
    "E2" = "-2.37%"
    "E3" = "-1.74%"
    "D4" = "5.061"
    "E4" = "-2.52%"
    "D5" = "0.07611"
    "E5" = "-4.92%"
    "D6" = "4.246"
    "E6" = "-2.99%"
    "D7" = "1.596"
    "E7" = "-7.03%"
    "E8" = "-1.96%"
    "D9" = "0.1009"
    "E9" = "-10.11%"
    "D10" = "0.1770"
    "E10" = "-3.48%"
    "D11" = "0.09140"
    "E11" = "0.02%"
    "D12" = "0.04386"
    "E12" = "-4.00%"
    "E13" = "-0.17%"
    "D14" = "0.001252"
    "E14" = "-2.19%"
    "D15" = "0.005863"
    "E15" = "-2.06%"
    "E16" = "0.33%"
    "D17" = "2.442"
    "E17" = "-4.97%"
    "E18" = "-2.50%"
    "D19" = "6.802"
    "E19" = "-7.43%"
    "D20" = "0.1357"
    "E20" = "-2.14%"
    "E21" = "7.69%"
    "D22" = "0.04156"
    "E22" = "-0.27%"
    "D23" = "0.001209"
    "E23" = "-3.42%"
    "D24" = "0.004064"
    "E24" = "-4.05%"
    "E25" = "5.46%"
    "D26" = "0.0003009"
    "E26" = "0.36%"
    "D38" = "0.02404"
    "E38" = "-5.18%"
    "D39" = "0.05143"
    "E39" = "-4.32%"
    "D40" = "0.007765"
    "E40" = "-4.07%"
    "E41" = "-5.39%"
    "D42" = "0.007091"
    "E42" = "-7.05%"
    "D43" = "0.001949"
    "E43" = "-6.35%"
    "D44" = "0.008387"
    "E44" = "0.33%"
    "D45" = "0.3054"
    "E45" = "-2.28%"
    "D46" = "0.00006370"
    "E46" = "-6.04%"
    "E47" = "-1.01%"
    "D48" = "0.006459"
    "E48" = "89.00%"
    "E49" = "-27.22%"
    "D50" = "0.00002101"
    "E50" = "-1.01%"
    "D51" = "0.0002001"
    "E51" = "-1.01%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}

